$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SLG building config update: after removing the Gamelogic project, most
# Func1..Func12 rows no longer trigger every EFT_* effect. Clear the
# corresponding flags (set to 0) while leaving the few still-applicable
# ones untouched.

# Row 2 (Func1): keep B,C,D,O = 1 ; clear E..N
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0

# Row 3 (Func2): keep B,E,O = 1 ; clear C,D,F..N
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0

# Row 4 (Func3): keep B,D,F,O = 1 ; clear C,E,G..N
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0

# Row 5 (Func4): keep B,E,O = 1 ; clear C,D,F..N
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0

# Rows 6-13 (Func5..Func12): keep B,O = 1 ; clear C..N
for ($r = 6; $r -le 13; $r++) {
    for ($c = 3; $c -le 14; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# Drop the unused built-in cell styles (thousands separator, currency,
# percent, etc.) that were never applied in this workbook - leaves only
# the default "Normal" style, matching the cleaned-up style table.
$styles = $wb.Styles
for ($i = $styles.Count; $i -ge 2; $i--) {
    $styles.Item($i).Delete()
}

# Update the last active selection recorded in the sheet view
$ws.Range("F11").Select()

$wb.Save()
